# Apply edits described by the diff:
# 1. Rename sheet from "o554F-HW30.xpc" to "o554F"
# 2. Tweak two floating point values in row 13 (G13, M13)
# 3. Add a new row 16 with data for "HexGrid-60degTilt5degRes" (reuses shared string index 13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet/tab
$ws.Name = "o554F"

# 2. Minor floating point corrections in existing row 13
$ws.Range("G13").Value = 0.9927666452891548
$ws.Range("M13").Value = 0.9937633963958241

# 3. Add new row 16, mirroring the structure of the other data rows
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.032673728103569
$ws.Range("D16").Value = 0.8643542507254005
$ws.Range("E16").Value = 1.021053435571271
$ws.Range("F16").Value = 1.032673728103569
$ws.Range("G16").Value = 0.9265996657552829
$ws.Range("H16").Value = 1.066101540864504
$ws.Range("I16").Value = 1.024103850120611
$ws.Range("J16").Value = 0.8643542507254005
$ws.Range("K16").Value = 0.9427038431483357
$ws.Range("L16").Value = 0.9876887856259524
$ws.Range("M16").Value = 0.9891477451901065

# Match the style of A column cells (bold, centered, bordered) used for A2:A15
$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").HorizontalAlignment = -4108
$ws.Range("A16").VerticalAlignment = -4160
$ws.Range("A16").Borders.LineStyle = 1
